# Append: 2025-09-13 06:30 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for the newly
# appended/refreshed rows (rows 2-7) on the active sheet ("ランサーズ")
# from "2025-09-13 06:23:43" to "2025-09-13 06:30:07".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-13 06:30:07"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
